$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the match data (columns F:V) between row 12 and row 13 ---
# Save row 12 values (F:V)
$r12 = @(
    $ws.Range("F12").Value(),
    $ws.Range("G12").Value(),
    $ws.Range("H12").Value(),
    $ws.Range("I12").Value(),
    $ws.Range("J12").Value(),
    $ws.Range("K12").Value(),
    $ws.Range("L12").Value(),
    $ws.Range("M12").Value(),
    $ws.Range("N12").Value(),
    $ws.Range("O12").Value(),
    $ws.Range("P12").Value(),
    $ws.Range("Q12").Value(),
    $ws.Range("R12").Value(),
    $ws.Range("S12").Value(),
    $ws.Range("T12").Value(),
    $ws.Range("U12").Value(),
    $ws.Range("V12").Value()
)

# Save row 13 values (F:V)
$r13 = @(
    $ws.Range("F13").Value(),
    $ws.Range("G13").Value(),
    $ws.Range("H13").Value(),
    $ws.Range("I13").Value(),
    $ws.Range("J13").Value(),
    $ws.Range("K13").Value(),
    $ws.Range("L13").Value(),
    $ws.Range("M13").Value(),
    $ws.Range("N13").Value(),
    $ws.Range("O13").Value(),
    $ws.Range("P13").Value(),
    $ws.Range("Q13").Value(),
    $ws.Range("R13").Value(),
    $ws.Range("S13").Value(),
    $ws.Range("T13").Value(),
    $ws.Range("U13").Value(),
    $ws.Range("V13").Value()
)

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $r13[$i]
    $ws.Range($cols[$i] + "13").Value = $r12[$i]
}

# --- 2) Append a new row 88 with the new match data ---
$ws.Range("A87:V87").Copy()
$ws.Range("A88:V88").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A88").Value = 87
$ws.Range("B88").Value = "south-africa"
$ws.Range("C88").Value = "premier-league"
$ws.Range("D88").Value = "2023-2024"
$ws.Range("E88").Value = 45242.69791666666
$ws.Range("F88").Value = "Sekhukhune"
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = "Supersport Utd"
$ws.Range("I88").Value = 1
$ws.Range("J88").Value = 3.59
$ws.Range("K88").Value = "08/11/2023 18:42"
$ws.Range("L88").Value = 3.51
$ws.Range("M88").Value = "12/11/2023 16:44"
$ws.Range("N88").Value = 2.9
$ws.Range("O88").Value = "08/11/2023 18:42"
$ws.Range("P88").Value = 2.83
$ws.Range("Q88").Value = "12/11/2023 16:44"
$ws.Range("R88").Value = 2.32
$ws.Range("S88").Value = "08/11/2023 18:42"
$ws.Range("T88").Value = 2.46
$ws.Range("U88").Value = "12/11/2023 16:44"
$ws.Range("V88").Value = "https://www.betexplorer.com/football/south-africa/premier-league/sekhukhune-supersport-utd/tjvrjKCE/"

Write-Host "Edit applied successfully"
